$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update volume/issue number and date-range headers
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# Update crime-statistics table (rows 14-30)
# Row 14
$ws.Range("N14").Value = -88.888888888888

# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -6.666666666666
$ws.Range("L15").Value = 27.272727272727
$ws.Range("M15").Value = 100

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -65
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = -38.541666666666
$ws.Range("L16").Value = -39.795918367346
$ws.Range("M16").Value = -44.859813084112
$ws.Range("N16").Value = -86.529680365296

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -12.121212121212
$ws.Range("I17").Value = 144
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 1.408450704225
$ws.Range("L17").Value = 10.769230769230
$ws.Range("M17").Value = 32.110091743119
$ws.Range("N17").Value = 11.627906976744

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("D18").Value = 2
$ws.Range("D18").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("E18").Value = -50
$ws.Range("E18").NumberFormat = $ws.Range("H18").NumberFormat
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = -11.363636363636
$ws.Range("L18").Value = 1.298701298701
$ws.Range("M18").Value = -57.377049180327
$ws.Range("N18").Value = -91.437980241492

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -43.75
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 12.244897959183
$ws.Range("I19").Value = 281
$ws.Range("J19").Value = 319
$ws.Range("K19").Value = -11.912225705329
$ws.Range("L19").Value = -5.704697986577
$ws.Range("M19").Value = 51.891891891891
$ws.Range("N19").Value = 3.690036900369

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -3.571428571428
$ws.Range("I20").Value = 134
$ws.Range("J20").Value = 162
$ws.Range("K20").Value = -17.283950617283
$ws.Range("L20").Value = -11.842105263157
$ws.Range("M20").Value = -18.292682926829
$ws.Range("N20").Value = -91.961607678464

# Row 21
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -38.888888888888
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = -12.587412587412
$ws.Range("I21").Value = 711
$ws.Range("J21").Value = 823
$ws.Range("K21").Value = -13.608748481166
$ws.Range("L21").Value = -7.421875
$ws.Range("M21").Value = -5.952380952380
$ws.Range("N21").Value = -79.325385286420

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = $ws.Range("H18").NumberFormat
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = $ws.Range("F18").NumberFormat
$ws.Range("H22").Value = 0
$ws.Range("H22").NumberFormat = $ws.Range("H18").NumberFormat
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -27.272727272727
$ws.Range("L22").Value = -38.461538461538
$ws.Range("M22").Value = -20

# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -10
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = -16.806722689075
$ws.Range("I24").Value = 551
$ws.Range("J24").Value = 644
$ws.Range("K24").Value = -14.440993788819
$ws.Range("L24").Value = -3.163444639718
$ws.Range("M24").Value = 11.764705882352

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = -18.181818181818
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 224
$ws.Range("J25").Value = 229
$ws.Range("K25").Value = -2.183406113537
$ws.Range("L25").Value = 22.404371584699

# Row 26
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -15.789473684210
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = -3.333333333333
$ws.Range("I26").Value = 245
$ws.Range("J26").Value = 280
$ws.Range("K26").Value = -12.5
$ws.Range("L26").Value = 6.521739130434
$ws.Range("M26").Value = -28.571428571428

# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -6.25

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 3.225806451612

# Row 29
$ws.Range("N29").Value = -83.333333333333

# Row 30
$ws.Range("N30").Value = -81.818181818181

